$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $value
}

function Set-PlainValue($row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.Value = $value
}

# Row 2
Set-TextValue 2 4 "256.42"
Set-TextValue 2 5 "-0.08%"

# Row 3
Set-TextValue 3 5 "-1.89%"

# Row 4
Set-TextValue 4 4 "4.649"
Set-TextValue 4 5 "-0.39%"

# Row 5
Set-TextValue 5 4 "0.05926"
Set-TextValue 5 5 "0.43%"

# Row 6
Set-TextValue 6 4 "6.606"
Set-TextValue 6 5 "-0.67%"

# Row 7
Set-TextValue 7 4 "0.8559"
Set-TextValue 7 5 "-1.37%"

# Row 8
Set-TextValue 8 4 "0.9113"
Set-TextValue 8 5 "-4.27%"

# Row 9
Set-TextValue 9 5 "-1.63%"

# Row 10
Set-TextValue 10 4 "0.04295"
Set-TextValue 10 5 "15.26%"

# Row 11
Set-TextValue 11 4 "0.07010"
Set-TextValue 11 5 "-0.98%"

# Row 12
Set-TextValue 12 4 "0.03023"
Set-TextValue 12 5 "-5.71%"

# Row 13
Set-TextValue 13 4 "0.09107"
Set-TextValue 13 5 "-1.66%"

# Row 14
Set-TextValue 14 4 "0.001527"
Set-TextValue 14 5 "-1.60%"

# Row 15
Set-PlainValue 15 2 "TigerCash"
Set-PlainValue 15 3 "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue 15 4 "0.006052"
Set-TextValue 15 5 "0.46%"

# Row 16
Set-PlainValue 16 2 "LEO"
Set-PlainValue 16 3 "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue 16 4 "3.475"
Set-TextValue 16 5 "-1.09%"

# Row 17
Set-PlainValue 17 2 "GateToken"
Set-PlainValue 17 3 "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue 17 4 "3.134"
Set-TextValue 17 5 "-1.78%"

# Row 18
Set-PlainValue 18 2 "BTSEToken"
Set-PlainValue 18 3 "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue 18 4 "2.150"
Set-TextValue 18 5 "-3.28%"

# Row 19
Set-PlainValue 19 2 "One"
Set-PlainValue 19 3 "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue 19 4 "0.01035"
Set-TextValue 19 5 "1,620.23%"

# Row 20
Set-TextValue 20 4 "0.3081"
Set-TextValue 20 5 "0.14%"

# Row 21
Set-TextValue 21 4 "0.1286"
Set-TextValue 21 5 "0.25%"

# Row 22
Set-TextValue 22 4 "3.890"
Set-TextValue 22 5 "1.14%"

# Row 24
Set-TextValue 24 4 "0.001216"
Set-TextValue 24 5 "-0.27%"

# Row 25
Set-TextValue 25 4 "0.004652"
Set-TextValue 25 5 "8.72%"

# Row 26
Set-TextValue 26 5 "0.01%"

# Row 27
Set-TextValue 27 4 "0.0001714"
Set-TextValue 27 5 "13.76%"

# Row 40
Set-TextValue 40 4 "0.03801"
Set-TextValue 40 5 "-0.43%"

# Row 41
Set-PlainValue 41 2 "KickToken"
Set-PlainValue 41 3 "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue 41 4 "0.006217"
Set-TextValue 41 5 "0.02%"

# Row 42
Set-PlainValue 42 2 "BKEXToken"
Set-PlainValue 42 3 "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue 42 4 "0.1099"
Set-TextValue 42 5 "0.08%"

# Row 43
Set-TextValue 43 4 "0.002199"
Set-TextValue 43 5 "-2.15%"

# Row 44
Set-TextValue 44 4 "0.01434"
Set-TextValue 44 5 "24.06%"

# Row 45
Set-TextValue 45 4 "0.00005136"
Set-TextValue 45 5 "-6.57%"

# Row 46
Set-TextValue 46 5 "0.01%"

# Row 47
Set-TextValue 47 4 "0.04998"
Set-TextValue 47 5 "-16.94%"

# Row 48
Set-TextValue 48 5 "10,470.24%"

# Row 49
Set-TextValue 49 4 "0.00002099"
Set-TextValue 49 5 "0.01%"

# Row 50
Set-TextValue 50 4 "0.0001999"
Set-TextValue 50 5 "0.01%"
